# Auto update Excel log
# Appends new sensor-event rows to the "Proximity", "mmWave" and "Camera"
# sheets of the SeniorConnect master log workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Proximity sheet: append rows 10-17 (ENTER/EXIT events for the
# Living Room Main Door).
# ---------------------------------------------------------------------
$wsProximity = $wb.Worksheets.Item("Proximity")

$proximityRows = @(
  @("2026-02-01","13:18:20","13:00","Living Room Main Door","ENTER","User ENTERED Living Room Main Door"),
  @("2026-02-01","13:18:50","13:00","Living Room Main Door","EXIT","User EXITED Living Room Main Door"),
  @("2026-02-01","13:18:56","13:00","Living Room Main Door","ENTER","User ENTERED Living Room Main Door"),
  @("2026-02-01","13:36:40","13:00","Living Room Main Door","ENTER","User ENTERED Living Room Main Door"),
  @("2026-02-01","13:36:43","13:00","Living Room Main Door","EXIT","User EXITED Living Room Main Door"),
  @("2026-02-01","13:36:55","13:00","Living Room Main Door","ENTER","User ENTERED Living Room Main Door"),
  @("2026-02-01","13:37:21","13:00","Living Room Main Door","EXIT","User EXITED Living Room Main Door"),
  @("2026-02-01","13:37:30","13:00","Living Room Main Door","ENTER","User ENTERED Living Room Main Door")
)

$startRow = 10
for ($i = 0; $i -lt $proximityRows.Count; $i++) {
  $r = $startRow + $i
  $values = $proximityRows[$i]
  for ($c = 1; $c -le 6; $c++) {
    $cell = $wsProximity.Cells.Item($r, $c)
    $cell.NumberFormat = "@"
    $cell.Value = $values[$c - 1]
  }
}

# ---------------------------------------------------------------------
# mmWave sheet: append rows 2-4 (PRESENCE_DETECTED events for the
# Living Room).
# ---------------------------------------------------------------------
$wsMmWave = $wb.Worksheets.Item("mmWave")

$mmWaveRows = @(
  @("2026-02-01","13:18:33","13:00","Living Room","PRESENCE_DETECTED","Active"),
  @("2026-02-01","13:18:50","13:00","Living Room","PRESENCE_DETECTED","Active"),
  @("2026-02-01","13:18:54","13:00","Living Room","PRESENCE_DETECTED","Active")
)

$startRow = 2
for ($i = 0; $i -lt $mmWaveRows.Count; $i++) {
  $r = $startRow + $i
  $values = $mmWaveRows[$i]
  for ($c = 1; $c -le 6; $c++) {
    $cell = $wsMmWave.Cells.Item($r, $c)
    $cell.NumberFormat = "@"
    $cell.Value = $values[$c - 1]
  }
}

# ---------------------------------------------------------------------
# Camera sheet: append rows 2-5 (Image Received events for the
# Living Room Main Door).
# ---------------------------------------------------------------------
$wsCamera = $wb.Worksheets.Item("Camera")

$cameraRows = @(
  @("2026-02-01","13:36:41","13:00","Living Room Main Door","Image Received","Active"),
  @("2026-02-01","13:36:56","13:00","Living Room Main Door","Image Received","Active"),
  @("2026-02-01","13:37:21","13:00","Living Room Main Door","Image Received","Active"),
  @("2026-02-01","13:37:32","13:00","Living Room Main Door","Image Received","Active")
)

$startRow = 2
for ($i = 0; $i -lt $cameraRows.Count; $i++) {
  $r = $startRow + $i
  $values = $cameraRows[$i]
  for ($c = 1; $c -le 6; $c++) {
    $cell = $wsCamera.Cells.Item($r, $c)
    $cell.NumberFormat = "@"
    $cell.Value = $values[$c - 1]
  }
}
